$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 95; existing rows 95-144 shift down to 96-145
$ws.Rows.Item(95).Insert()

# Populate the new row 95 with the new weekly record
$ws.Cells.Item(95, 1).Value = 7
$ws.Cells.Item(95, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(95, 3).Value = "Ñuble"
$ws.Cells.Item(95, 4).Value = 44455
$ws.Cells.Item(95, 4).NumberFormat = $ws.Cells.Item(96, 4).NumberFormat
$ws.Cells.Item(95, 5).Value = 16
$ws.Cells.Item(95, 6).Value = 100112043
$ws.Cells.Item(95, 7).Value = "Pepino ensalada"
$ws.Cells.Item(95, 8).Value = "Sin especificar"
$ws.Cells.Item(95, 9).Value = "Primera"
$ws.Cells.Item(95, 10).Value = 300
$ws.Cells.Item(95, 11).Value = 16000
$ws.Cells.Item(95, 12).Value = 17000
$ws.Cells.Item(95, 13).Value = 16500
$ws.Cells.Item(95, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(95, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(95, 16).Value = 275
$ws.Cells.Item(95, 17).Value = 60
$ws.Cells.Item(95, 18).Value = "Hortaliza"
